$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 590
$ws.Range("I55").Value = 533.3333
$ws.Range("J55").Value = 675
$ws.Range("K55").Value = 533.3333
$ws.Range("L55").Value = 675
$ws.Range("M55").Value = -319.3333
$ws.Range("N55").Value = -1103
$ws.Range("H112").Value = 2270.3845
$ws.Range("J112").Value = 2376.5
$ws.Range("L112").Value = 7129.5
$ws.Range("N112").Value = -9345.5
$ws.Range("H132").Value = 2666.6438
$ws.Range("I132").Value = 2691.3088
$ws.Range("J132").Value = 2331.2
$ws.Range("K132").Value = 8073.926399999999
$ws.Range("L132").Value = 6993.599999999999
$ws.Range("M132").Value = -5543.926399999999
$ws.Range("N132").Value = -12053.6
$ws.Range("H135").Value = 1178.8158
$ws.Range("I135").Value = 1178.8158
$ws.Range("K135").Value = 10609.3422
$ws.Range("M135").Value = -8074.342200000001
$ws.Range("H137").Value = 3206
$ws.Range("I137").Value = 2400
$ws.Range("J137").Value = 3279.2727
$ws.Range("K137").Value = 7200
$ws.Range("L137").Value = 9837.8181
$ws.Range("M137").Value = -4650
$ws.Range("N137").Value = -14937.8181
$ws.Range("H138").Value = 3585.8064
$ws.Range("I138").Value = 3245.4614
$ws.Range("J138").Value = 3831.611
$ws.Range("K138").Value = 9736.3842
$ws.Range("L138").Value = 11494.833
$ws.Range("M138").Value = -4596.3842
$ws.Range("N138").Value = -21774.833

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1604.3467
$ws.Range("I32").Value = 839.65717
$ws.Range("J32").Value = 12310
$ws.Range("K32").Value = 839.65717
$ws.Range("L32").Value = 12310
$ws.Range("M32").Value = -552.65717
$ws.Range("N32").Value = -12884
$ws.Range("H61").Value = 2473.8333
$ws.Range("I61").Value = 2473.8333
$ws.Range("K61").Value = 2473.8333
$ws.Range("M61").Value = -2261.8333
$ws.Range("H74").Value = 1952.7174
$ws.Range("I74").Value = 1201.7368
$ws.Range("J74").Value = 2481.1853
$ws.Range("K74").Value = 1201.7368
$ws.Range("L74").Value = 2481.1853
$ws.Range("M74").Value = -327.7367999999999
$ws.Range("N74").Value = -4229.1853
$ws.Range("H77").Value = 1952.7174
$ws.Range("I77").Value = 1201.7368
$ws.Range("J77").Value = 2481.1853
$ws.Range("K77").Value = 6008.683999999999
$ws.Range("L77").Value = 12405.9265
$ws.Range("M77").Value = -1640.683999999999
$ws.Range("N77").Value = -21141.9265
$ws.Range("H132").Value = 5342.1665
$ws.Range("I132").Value = 4697.3335
$ws.Range("J132").Value = 5987
$ws.Range("K132").Value = 14092.0005
$ws.Range("L132").Value = 17961
$ws.Range("M132").Value = -11562.0005
$ws.Range("N132").Value = -23021
$ws.Range("H136").Value = 2473.8333
$ws.Range("I136").Value = 2473.8333
$ws.Range("K136").Value = 7421.499899999999
$ws.Range("M136").Value = -4871.499899999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 574.7727
$ws.Range("I22").Value = 379.88235
$ws.Range("J22").Value = 1237.4
$ws.Range("K22").Value = 379.88235
$ws.Range("L22").Value = 1237.4
$ws.Range("M22").Value = -206.88235
$ws.Range("N22").Value = -1583.4
$ws.Range("H134").Value = 4781.1577
$ws.Range("I134").Value = 4037.5
$ws.Range("J134").Value = 7298.154
$ws.Range("K134").Value = 12112.5
$ws.Range("L134").Value = 21894.462
$ws.Range("M134").Value = -9577.5
$ws.Range("N134").Value = -26964.462

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3188.5925
$ws.Range("I31").Value = 2448.087
$ws.Range("K31").Value = 2448.087
$ws.Range("M31").Value = -2153.087
$ws.Range("H34").Value = 3188.5925
$ws.Range("I34").Value = 2448.087
$ws.Range("K34").Value = 2448.087
$ws.Range("M34").Value = -2246.087
$ws.Range("H58").Value = 1638.1471
$ws.Range("I58").Value = 1145
$ws.Range("J58").Value = 2669.2727
$ws.Range("K58").Value = 1145
$ws.Range("L58").Value = 2669.2727
$ws.Range("M58").Value = -942
$ws.Range("N58").Value = -3075.2727
$ws.Range("H132").Value = 3725.5122
$ws.Range("I132").Value = 3398.5715
$ws.Range("J132").Value = 5632.6665
$ws.Range("K132").Value = 10195.7145
$ws.Range("L132").Value = 16897.9995
$ws.Range("M132").Value = -7665.7145
$ws.Range("N132").Value = -21957.9995
$ws.Range("H134").Value = 3960.1482
$ws.Range("I134").Value = 3477
$ws.Range("K134").Value = 10431
$ws.Range("M134").Value = -7896
$ws.Range("H136").Value = 1638.1471
$ws.Range("I136").Value = 1145
$ws.Range("J136").Value = 2669.2727
$ws.Range("K136").Value = 3435
$ws.Range("L136").Value = 8007.8181
$ws.Range("M136").Value = -885
$ws.Range("N136").Value = -13107.8181
$ws.Range("H141").Value = 273333
$ws.Range("J141").Value = 273333
$ws.Range("L141").Value = 273333
$ws.Range("N141").Value = -283693

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 1451.3846
$ws.Range("I113").Value = 1344.5
$ws.Range("J113").Value = 1498.8889
$ws.Range("K113").Value = 4033.5
$ws.Range("L113").Value = 4496.6667
$ws.Range("M113").Value = -1863.5
$ws.Range("N113").Value = -8836.6667
$ws.Range("H131").Value = 17095572
$ws.Range("I131").Value = 10102340
$ws.Range("J131").Value = 22223944
$ws.Range("K131").Value = 30307020
$ws.Range("L131").Value = 66671832
$ws.Range("M131").Value = -30301980
$ws.Range("N131").Value = -66681912
$ws.Range("H134").Value = 17386.666
$ws.Range("I134").Value = 5800
$ws.Range("K134").Value = 17400
$ws.Range("M134").Value = -12330
$ws.Range("H137").Value = 7786.7144
$ws.Range("I137").Value = 2756
$ws.Range("J137").Value = 14494.333
$ws.Range("K137").Value = 8268
$ws.Range("L137").Value = 43482.999
$ws.Range("M137").Value = -3168
$ws.Range("N137").Value = -53682.999
$ws.Range("H139").Value = 11119241
$ws.Range("I139").Value = 12825278
$ws.Range("J139").Value = 30000
$ws.Range("K139").Value = 38475834
$ws.Range("L139").Value = 90000
$ws.Range("M139").Value = -38470694
$ws.Range("N139").Value = -100280
$ws.Range("H140").Value = 5329745
$ws.Range("I140").Value = 25002788
$ws.Range("J140").Value = 12706.703
$ws.Range("K140").Value = 75008364
$ws.Range("L140").Value = 38120.109
$ws.Range("M140").Value = -75003184
$ws.Range("N140").Value = -48480.109
$ws.Range("H141").Value = 24744.426
$ws.Range("I141").Value = 6370.3335
$ws.Range("K141").Value = 19111.0005
$ws.Range("M141").Value = -13931.0005

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2882.7144
$ws.Range("I80").Value = 2794.75
$ws.Range("K80").Value = 2794.75
$ws.Range("M80").Value = -1796.75
$ws.Range("H83").Value = 2882.7144
$ws.Range("I83").Value = 2794.75
$ws.Range("K83").Value = 13973.75
$ws.Range("M83").Value = -8981.75
$ws.Range("H126").Value = 7105.2
$ws.Range("I126").Value = 6560.4
$ws.Range("J126").Value = 7650
$ws.Range("K126").Value = 19681.2
$ws.Range("L126").Value = 22950
$ws.Range("M126").Value = -17211.2
$ws.Range("N126").Value = -27890
$ws.Range("H132").Value = 4552.25
$ws.Range("I132").Value = 5069.6665
$ws.Range("J132").Value = 3000
$ws.Range("K132").Value = 15208.9995
$ws.Range("L132").Value = 9000
$ws.Range("M132").Value = -12678.9995
$ws.Range("N132").Value = -14060

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1632.2354
$ws.Range("I22").Value = 1249.3334
$ws.Range("J22").Value = 1714.2858
$ws.Range("K22").Value = 1249.3334
$ws.Range("L22").Value = 1714.2858
$ws.Range("M22").Value = -954.3334
$ws.Range("N22").Value = -2304.2858
$ws.Range("H27").Value = 1632.2354
$ws.Range("I27").Value = 1249.3334
$ws.Range("J27").Value = 1714.2858
$ws.Range("K27").Value = 1249.3334
$ws.Range("L27").Value = 1714.2858
$ws.Range("M27").Value = -1142.3334
$ws.Range("N27").Value = -1928.2858
$ws.Range("H132").Value = 2404.3242
$ws.Range("I132").Value = 2062.3809
$ws.Range("K132").Value = 6187.1427
$ws.Range("M132").Value = -3657.1427

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2840.5
$ws.Range("I132").Value = 2815.6592
$ws.Range("J132").Value = 2977.125
$ws.Range("K132").Value = 8446.9776
$ws.Range("L132").Value = 8931.375
$ws.Range("M132").Value = -5916.9776
$ws.Range("N132").Value = -13991.375
$ws.Range("H136").Value = 2240.1765
$ws.Range("I136").Value = 2203.5186
$ws.Range("K136").Value = 6610.5558
$ws.Range("M136").Value = -4060.5558
